$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting the existing rows 35-44 down to 36-45
$ws.Rows("35").Insert()

# Fill in the new weekly price record for row 35
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C35").Value = 'Los Lagos'
$ws.Range("D35").Value = 44663
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 100112030
$ws.Range("G35").Value = 'Poroto granado'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 60
$ws.Range("K35").Value = 28000
$ws.Range("L35").Value = 28000
$ws.Range("M35").Value = 28000
$ws.Range("N35").Value = '$/saco 25 kilos'
$ws.Range("O35").Value = 'Región Metropolitana'
$ws.Range("P35").Value = 1120
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = 'Hortaliza'
